$p = $ppt.ActivePresentation

# --- Slide 1: "TextBox 3" shape (title "2nd  프로젝트 개요 발표") ---
$s1 = $p.Slides.Item(1)
$sh1 = $s1.Shapes.Item(1)

# Reposition / resize the shape (Shape.Left/Width are in points; the
# literals below are the nearest double to 3174829/12700 EMU and
# 5841663/12700 EMU whose round-trip through PowerPoint's internal
# single-precision (float32) storage lands back on the exact target
# EMU value instead of the adjacent one).
$sh1.Left = 249.98654174804688
$sh1.Width = 459.9734802246094

# Split the single run "2nd  " into "2" (sz44) + "nd" (shrunk to sz28)
# + "  " (sz44), matching the superscript-style "nd" styling.
$tr1 = $sh1.TextFrame.TextRange
$tr1.Characters(2, 2).Font.Size = 28

# --- Slide 2: "TextBox 5" shape (bullet list incl. "일 평균 주행거리가 매년 줄어 들고 있음") ---
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(4)
$tr2 = $sh2.TextFrame.TextRange

# Merge the three adjacent runs "일 " + "평균 주행거리가 " + "매년 줄어 들고 있음"
# (characters 97-118 of the text frame) into a single run.
$tr2.Characters(97, 22).Text = "일 평균 주행거리가 매년 줄어 들고 있음"

# The shape auto-fits its text (<a:spAutoFit/>); shortening the text made
# PowerPoint shrink the box height automatically. The diff doesn't touch
# this shape's size/position, so restore the original height (the literal
# is the nearest double to 2862322/12700 EMU that round-trips exactly).
$sh2.Height = 225.37969970703125
